# Weekly update: insert a new price-report row for Coliflor at
# Feria Lagunitas de Puerto Montt, right above the existing row 194.
# All rows from 194..228 shift down by one (to 195..229); dimension
# grows from A1:R228 to A1:R229.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 194, pushing rows 194-228 down to 195-229.
$ws.Rows.Item(194).Insert()

# Populate the newly inserted row 194 with the new weekly data point.
$ws.Range("A194").Value2 = 4
$ws.Range("B194").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C194").Value2 = "Los Lagos"
$ws.Range("D194").Value2 = 44522
$ws.Range("E194").Value2 = 10
$ws.Range("F194").Value2 = 100112008
$ws.Range("G194").Value2 = "Coliflor"
$ws.Range("H194").Value2 = "Sin especificar"
$ws.Range("I194").Value2 = "Segunda"
$ws.Range("J194").Value2 = 500
$ws.Range("K194").Value2 = 1000
$ws.Range("L194").Value2 = 1000
$ws.Range("M194").Value2 = 1000
$ws.Range("N194").Value2 = "$/unidad"
$ws.Range("O194").Value2 = "Región Metropolitana"
$ws.Range("P194").Value2 = 1000
$ws.Range("Q194").Value2 = 1
$ws.Range("R194").Value2 = "Hortaliza"
